$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be explicit Text (t="str"/"s"), even for "numeric-looking"
# or empty strings, by writing it through the classic Lotus-style text-prefix apostrophe
# (Formula = "'" + value) rather than .Value, which otherwise tries to infer a numeric
# type, or drops an empty string entirely instead of keeping an empty Text cell.
function Set-TextCell($rng, $value) {
    $rng.Formula = "'" + $value
}

# Row 8
$ws.Range("C8").Value = 56
Set-TextCell $ws.Range("D8") '1.0'
Set-TextCell $ws.Range("E8") 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'

# Row 9
Set-TextCell $ws.Range("A9") 'Each'
$ws.Range("C9").Value = 41
Set-TextCell $ws.Range("D9") '4.0'
Set-TextCell $ws.Range("E9") 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F9").Value = 50
Set-TextCell $ws.Range("G9") '2050.00'

# Row 10
Set-TextCell $ws.Range("A10") 'R. mtr.'
$ws.Range("C10").Value = 84
Set-TextCell $ws.Range("D10") '17'
Set-TextCell $ws.Range("E10") '25 mm'
$ws.Range("F10").Value = 56
Set-TextCell $ws.Range("G10") '4704.00'

# Row 11
Set-TextCell $ws.Range("A11") 'Set'
$ws.Range("C11").Value = 58
Set-TextCell $ws.Range("D11") '13.0'
Set-TextCell $ws.Range("E11") 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F11").Value = 5733
Set-TextCell $ws.Range("G11") '332514.00'

# Row 12
Set-TextCell $ws.Range("A12") ''
$ws.Range("C12").Value = 2
Set-TextCell $ws.Range("D12") '15.0'
Set-TextCell $ws.Range("E12") 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 0
Set-TextCell $ws.Range("G12") '0.00'

# Row 13
Set-TextCell $ws.Range("A13") ''
$ws.Range("C13").Value = 51
Set-TextCell $ws.Range("D13") '16.0'
Set-TextCell $ws.Range("E13") 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 0
Set-TextCell $ws.Range("G13") '0.00'

# Row 14
$ws.Range("C14").Value = 75
Set-TextCell $ws.Range("D14") '31'
Set-TextCell $ws.Range("E14") 'Double pole MCB(With B/C curve tripping Characteristics)'

# Row 15
Set-TextCell $ws.Range("A15") 'Each'
$ws.Range("C15").Value = 85
Set-TextCell $ws.Range("D15") '35'
Set-TextCell $ws.Range("E15") '8 Way (8+2)'
$ws.Range("F15").Value = 2184
Set-TextCell $ws.Range("G15") '185640.00'

# Row 16
Set-TextCell $ws.Range("A16") '%'
$ws.Range("C16").Value = 51
Set-TextCell $ws.Range("D16") '37'
Set-TextCell $ws.Range("E16") 'Add Tender Premium '

# Remove the stale "Grand Total" row; this shifts rows 18-21 up to 17-20,
# which already carries the right labels (Grand Total Rs. / Tender Premium / NET PAYABLE)
# one row earlier - matching the target layout.
$ws.Rows(17).Delete()

# Update the totals to the new Grand Total value
Set-TextCell $ws.Range("G18") '524908.00'
Set-TextCell $ws.Range("H18") '524908.00'
Set-TextCell $ws.Range("G20") '524908.00'
Set-TextCell $ws.Range("H20") '524908.00'

